$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M34").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("H58").Value = 1016.1429
$ws.Range("I58").Value = 222.6
$ws.Range("K58").Value = 667.8
$ws.Range("M58").Value = -517.8
$ws.Range("H62").Value = 3950.5625
$ws.Range("I62").Value = 3960.75
$ws.Range("J62").Value = 3920
$ws.Range("K62").Value = 3960.75
$ws.Range("L62").Value = 3920
$ws.Range("M62").Value = -3336.75
$ws.Range("N62").Value = -5168
$ws.Range("H65").Value = 3950.5625
$ws.Range("I65").Value = 3960.75
$ws.Range("J65").Value = 3920
$ws.Range("K65").Value = 19803.75
$ws.Range("L65").Value = 19600
$ws.Range("M65").Value = -16683.75
$ws.Range("N65").Value = -25840
$ws.Range("H107").Value = 583.5
$ws.Range("I107").Value = 124
$ws.Range("K107").Value = 124
$ws.Range("M107").Value = 1796
$ws.Range("H115").Value = 359.33334
$ws.Range("I115").Value = 359.33334
$ws.Range("K115").Value = 1078.00002
$ws.Range("M115").Value = 488.9999800000001
$ws.Range("H118").Value = 341.8889
$ws.Range("I118").Value = 341.8889
$ws.Range("K118").Value = 1025.6667
$ws.Range("M118").Value = 631.3333
$ws.Range("H132").Value = 2140.4482
$ws.Range("I132").Value = 2038.3214
$ws.Range("K132").Value = 6114.9642
$ws.Range("M132").Value = -3584.9642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 283.16666
$ws.Range("I5").Value = 279.8
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 279.8
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -167.8
$ws.Range("N5").Value = -524
$ws.Range("N92").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 283.16666
$ws.Range("I4").Value = 279.8
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 279.8
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -164.8
$ws.Range("N4").Value = -530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 418.85715
$ws.Range("I22").Value = 350
$ws.Range("K22").Value = 350
$ws.Range("M22").Value = 0
$ws.Range("H122").Value = 900.2222
$ws.Range("I122").Value = 886
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 2658
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -208
$ws.Range("N122").Value = -7942
$ws.Range("H132").Value = 1090.4166
$ws.Range("I132").Value = 871.36365
$ws.Range("K132").Value = 2614.09095
$ws.Range("M132").Value = -84.09094999999979
$ws.Range("H134").Value = 2584.0557
$ws.Range("I134").Value = 2275.75
$ws.Range("K134").Value = 6827.25
$ws.Range("M134").Value = -4292.25
$ws.Range("N140").ClearContents()
$ws.Range("H140").Value = 124999
$ws.Range("I140").Value = 124999
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 124999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -119819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 310.91666
$ws.Range("I7").Value = 226.66667
$ws.Range("K7").Value = 680.00001
$ws.Range("M7").Value = -568.00001
$ws.Range("H17").Value = 1009
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 1332
$ws.Range("K17").Value = 120
$ws.Range("L17").Value = 3996
$ws.Range("M17").Value = 49
$ws.Range("N17").Value = -4334
$ws.Range("H26").Value = 105723.1
$ws.Range("J26").Value = 10284.2
$ws.Range("L26").Value = 30852.6
$ws.Range("N26").Value = -31428.6
$ws.Range("H68").Value = 1783.3334
$ws.Range("I68").Value = 1850
$ws.Range("J68").Value = 1750
$ws.Range("K68").Value = 5550
$ws.Range("L68").Value = 5250
$ws.Range("M68").Value = -4739
$ws.Range("N68").Value = -6872
$ws.Range("H71").Value = 1783.3334
$ws.Range("I71").Value = 1850
$ws.Range("J71").Value = 1750
$ws.Range("K71").Value = 16650
$ws.Range("L71").Value = 15750
$ws.Range("M71").Value = -12594
$ws.Range("N71").Value = -23862
$ws.Range("H107").Value = 839.17645
$ws.Range("J107").Value = 770.63635
$ws.Range("L107").Value = 2311.90905
$ws.Range("N107").Value = -6151.90905
$ws.Range("H121").Value = 839.6667
$ws.Range("I121").Value = 773.2857
$ws.Range("J121").Value = 897.75
$ws.Range("K121").Value = 2319.8571
$ws.Range("L121").Value = 2693.25
$ws.Range("M121").Value = -1009.8571
$ws.Range("N121").Value = -5313.25
$ws.Range("H140").Value = 1526.1111
$ws.Range("I140").Value = 1526.1111
$ws.Range("K140").Value = 4578.3333
$ws.Range("M140").Value = 601.6666999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M20").ClearContents()
$ws.Range("H20").Value = 52000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 52000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 52000
$ws.Range("N20").Value = -52490
$ws.Range("H24").Value = 35006
$ws.Range("I24").Value = 35006
$ws.Range("K24").Value = 35006
$ws.Range("M24").Value = -34833
$ws.Range("H80").Value = 2827.1
$ws.Range("I80").Value = 2469.1428
$ws.Range("J80").Value = 3662.3333
$ws.Range("K80").Value = 2469.1428
$ws.Range("L80").Value = 3662.3333
$ws.Range("M80").Value = -1471.1428
$ws.Range("N80").Value = -5658.3333
$ws.Range("H83").Value = 2827.1
$ws.Range("I83").Value = 2469.1428
$ws.Range("J83").Value = 3662.3333
$ws.Range("K83").Value = 12345.714
$ws.Range("L83").Value = 18311.6665
$ws.Range("M83").Value = -7353.714
$ws.Range("N83").Value = -28295.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 6500
$ws.Range("J20").Value = 6333.3335
$ws.Range("L20").Value = 6333.3335
$ws.Range("N20").Value = -6785.3335
$ws.Range("H21").Value = 11000
$ws.Range("I21").Value = 14000
$ws.Range("J21").Value = 8000
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = -13826
$ws.Range("N21").Value = -8348
$ws.Range("N43").ClearContents()
$ws.Range("I43").Value = 12000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 12000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -11807

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1249.75
$ws.Range("I126").Value = 999.5
$ws.Range("K126").Value = 2998.5
$ws.Range("M126").Value = -528.5
